# Auto-generated Excel COM-interop script
# Applies the diff: bumps C column (Förändrad) to 46070 for all data rows,
# and re-syncs rows 5-41 (Beteckning/Datum/Area + the two 'Kommuner' Markägare tags)
# to match the refreshed source extract.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A 22581-2025 -> A 22581-2025
$ws.Cells.Item(2, 3).Value = 46070

# Row 3: A 56290-2025 -> A 56290-2025
$ws.Cells.Item(3, 3).Value = 46070

# Row 4: A 889-2026 -> A 889-2026
$ws.Cells.Item(4, 3).Value = 46070

# Row 5: A 59011-2025 -> A 23295-2025
$ws.Cells.Item(5, 1).Value = "A 23295-2025"
$ws.Cells.Item(5, 2).Value = 45791.59071759259
$ws.Cells.Item(5, 3).Value = 46070
$ws.Cells.Item(5, 7).Value = 2.2

# Row 6: A 1103-2025 -> A 23294-2025
$ws.Cells.Item(6, 1).Value = "A 23294-2025"
$ws.Cells.Item(6, 2).Value = 45791.58967592593
$ws.Cells.Item(6, 3).Value = 46070
$ws.Cells.Item(6, 7).Value = 4.6

# Row 7: A 21141-2023 -> A 23301-2025
$ws.Cells.Item(7, 1).Value = "A 23301-2025"
$ws.Cells.Item(7, 2).Value = 45791.59498842592
$ws.Cells.Item(7, 3).Value = 46070
$ws.Cells.Item(7, 7).Value = 2.9

# Row 8: A 34508-2025 -> A 42462-2025
$ws.Cells.Item(8, 1).Value = "A 42462-2025"
$ws.Cells.Item(8, 2).Value = 45905.45730324074
$ws.Cells.Item(8, 3).Value = 46070
$ws.Cells.Item(8, 7).Value = 2.3

# Row 9: A 34394-2025 -> A 1103-2025
$ws.Cells.Item(9, 1).Value = "A 1103-2025"
$ws.Cells.Item(9, 2).Value = 45666
$ws.Cells.Item(9, 3).Value = 46070
$ws.Cells.Item(9, 7).Value = 1.6

# Row 10: A 32577-2025 -> A 45167-2025
$ws.Cells.Item(10, 1).Value = "A 45167-2025"
$ws.Cells.Item(10, 2).Value = 45919.49364583333
$ws.Cells.Item(10, 3).Value = 46070
$ws.Cells.Item(10, 7).Value = 0.6

# Row 11: A 42462-2025 -> A 45088-2025
$ws.Cells.Item(11, 1).Value = "A 45088-2025"
$ws.Cells.Item(11, 2).Value = 45919.37598379629
$ws.Cells.Item(11, 3).Value = 46070
$ws.Cells.Item(11, 7).Value = 0.9

# Row 12: A 34939-2025 -> A 45158-2025
$ws.Cells.Item(12, 1).Value = "A 45158-2025"
$ws.Cells.Item(12, 2).Value = 45919.48245370371
$ws.Cells.Item(12, 3).Value = 46070
$ws.Cells.Item(12, 7).Value = 2.7

# Row 13: A 34963-2025 -> A 34508-2025
$ws.Cells.Item(13, 1).Value = "A 34508-2025"
$ws.Cells.Item(13, 2).Value = 45847.44315972222
$ws.Cells.Item(13, 3).Value = 46070
$ws.Cells.Item(13, 7).Value = 2

# Row 14: A 34591-2022 -> A 34394-2025
$ws.Cells.Item(14, 1).Value = "A 34394-2025"
$ws.Cells.Item(14, 2).Value = 45846.58854166666
$ws.Cells.Item(14, 3).Value = 46070
$ws.Cells.Item(14, 7).Value = 3.9

# Row 15: A 45167-2025 -> A 21141-2023
$ws.Cells.Item(15, 1).Value = "A 21141-2023"
$ws.Cells.Item(15, 2).Value = 45062
$ws.Cells.Item(15, 3).Value = 46070
$ws.Cells.Item(15, 7).Value = 3.4

# Row 16: A 45088-2025 -> A 32577-2025
$ws.Cells.Item(16, 1).Value = "A 32577-2025"
$ws.Cells.Item(16, 2).Value = 45838
$ws.Cells.Item(16, 3).Value = 46070
$ws.Cells.Item(16, 7).Value = 5.3

# Row 17: A 45158-2025 -> A 1379-2024
$ws.Cells.Item(17, 1).Value = "A 1379-2024"
$ws.Cells.Item(17, 2).Value = 45303.55193287037
$ws.Cells.Item(17, 3).Value = 46070
$ws.Cells.Item(17, 6).Value = "Kommuner"
$ws.Cells.Item(17, 7).Value = 1.5

# Row 18: A 60718-2022 -> A 61627-2024
$ws.Cells.Item(18, 1).Value = "A 61627-2024"
$ws.Cells.Item(18, 2).Value = 45646.66263888889
$ws.Cells.Item(18, 3).Value = 46070
$ws.Cells.Item(18, 6).Value = "Kommuner"
$ws.Cells.Item(18, 7).Value = 2.1

# Row 19: A 21379-2023 -> A 34939-2025
$ws.Cells.Item(19, 1).Value = "A 34939-2025"
$ws.Cells.Item(19, 2).Value = 45849.58229166667
$ws.Cells.Item(19, 3).Value = 46070
$ws.Cells.Item(19, 7).Value = 7.8

# Row 20: A 1379-2024 -> A 34963-2025
$ws.Cells.Item(20, 1).Value = "A 34963-2025"
$ws.Cells.Item(20, 2).Value = 45849.63219907408
$ws.Cells.Item(20, 3).Value = 46070
$ws.Cells.Item(20, 6).ClearContents()
$ws.Cells.Item(20, 7).Value = 1.1

# Row 21: A 61627-2024 -> A 34591-2022
$ws.Cells.Item(21, 1).Value = "A 34591-2022"
$ws.Cells.Item(21, 2).Value = 44795.3778587963
$ws.Cells.Item(21, 3).Value = 46070
$ws.Cells.Item(21, 6).ClearContents()
$ws.Cells.Item(21, 7).Value = 2.5

# Row 22: A 60392-2022 -> A 59011-2025
$ws.Cells.Item(22, 1).Value = "A 59011-2025"
$ws.Cells.Item(22, 2).Value = 45987
$ws.Cells.Item(22, 3).Value = 46070
$ws.Cells.Item(22, 7).Value = 2.7

# Row 23: A 7787-2023 -> A 60718-2022
$ws.Cells.Item(23, 1).Value = "A 60718-2022"
$ws.Cells.Item(23, 2).Value = 44912.89109953704
$ws.Cells.Item(23, 3).Value = 46070
$ws.Cells.Item(23, 7).Value = 1

# Row 24: A 60717-2022 -> A 57000-2025
$ws.Cells.Item(24, 1).Value = "A 57000-2025"
$ws.Cells.Item(24, 2).Value = 45977
$ws.Cells.Item(24, 3).Value = 46070
$ws.Cells.Item(24, 7).Value = 2.3

# Row 25: A 57000-2025 -> A 57655-2025
$ws.Cells.Item(25, 1).Value = "A 57655-2025"
$ws.Cells.Item(25, 2).Value = 45981.40369212963
$ws.Cells.Item(25, 3).Value = 46070
$ws.Cells.Item(25, 7).Value = 3.4

# Row 26: A 57655-2025 -> A 21379-2023
$ws.Cells.Item(26, 1).Value = "A 21379-2023"
$ws.Cells.Item(26, 2).Value = 45063.34819444444
$ws.Cells.Item(26, 3).Value = 46070
$ws.Cells.Item(26, 7).Value = 5.8

# Row 27: A 43229-2024 -> A 2434-2026
$ws.Cells.Item(27, 1).Value = "A 2434-2026"
$ws.Cells.Item(27, 2).Value = 46036.86722222222
$ws.Cells.Item(27, 3).Value = 46070
$ws.Cells.Item(27, 7).Value = 1.2

# Row 28: A 59192-2022 -> A 2433-2026
$ws.Cells.Item(28, 1).Value = "A 2433-2026"
$ws.Cells.Item(28, 2).Value = 46036.86631944445
$ws.Cells.Item(28, 3).Value = 46070
$ws.Cells.Item(28, 7).Value = 2.6

# Row 29: A 2434-2026 -> A 60392-2022
$ws.Cells.Item(29, 1).Value = "A 60392-2022"
$ws.Cells.Item(29, 2).Value = 44910
$ws.Cells.Item(29, 3).Value = 46070
$ws.Cells.Item(29, 7).Value = 4.1

# Row 30: A 7791-2023 -> A 7787-2023
$ws.Cells.Item(30, 1).Value = "A 7787-2023"
$ws.Cells.Item(30, 3).Value = 46070
$ws.Cells.Item(30, 7).Value = 5.7

# Row 31: A 2433-2026 -> A 60717-2022
$ws.Cells.Item(31, 1).Value = "A 60717-2022"
$ws.Cells.Item(31, 2).Value = 44912.89078703704
$ws.Cells.Item(31, 3).Value = 46070
$ws.Cells.Item(31, 7).Value = 0.7

# Row 32: A 769-2023 -> A 43229-2024
$ws.Cells.Item(32, 1).Value = "A 43229-2024"
$ws.Cells.Item(32, 2).Value = 45567.88697916667
$ws.Cells.Item(32, 3).Value = 46070
$ws.Cells.Item(32, 7).Value = 9.699999999999999

# Row 33: A 1486-2023 -> A 59192-2022
$ws.Cells.Item(33, 1).Value = "A 59192-2022"
$ws.Cells.Item(33, 2).Value = 44896
$ws.Cells.Item(33, 3).Value = 46070
$ws.Cells.Item(33, 7).Value = 1.5

# Row 34: A 60793-2023 -> A 7791-2023
$ws.Cells.Item(34, 1).Value = "A 7791-2023"
$ws.Cells.Item(34, 2).Value = 44973
$ws.Cells.Item(34, 3).Value = 46070
$ws.Cells.Item(34, 7).Value = 3.1

# Row 35: A 57410-2022 -> A 769-2023
$ws.Cells.Item(35, 1).Value = "A 769-2023"
$ws.Cells.Item(35, 2).Value = 44931
$ws.Cells.Item(35, 3).Value = 46070
$ws.Cells.Item(35, 7).Value = 1.7

# Row 36: A 31486-2021 -> A 1486-2023
$ws.Cells.Item(36, 1).Value = "A 1486-2023"
$ws.Cells.Item(36, 2).Value = 44937
$ws.Cells.Item(36, 3).Value = 46070
$ws.Cells.Item(36, 7).Value = 0.2

# Row 37: A 63664-2023 -> A 60793-2023
$ws.Cells.Item(37, 1).Value = "A 60793-2023"
$ws.Cells.Item(37, 2).Value = 45260.6534375
$ws.Cells.Item(37, 3).Value = 46070
$ws.Cells.Item(37, 7).Value = 0.7

# Row 38: A 57955-2024 -> A 57410-2022
$ws.Cells.Item(38, 1).Value = "A 57410-2022"
$ws.Cells.Item(38, 2).Value = 44896
$ws.Cells.Item(38, 3).Value = 46070
$ws.Cells.Item(38, 7).Value = 7.5

# Row 39: A 23295-2025 -> A 31486-2021
$ws.Cells.Item(39, 1).Value = "A 31486-2021"
$ws.Cells.Item(39, 2).Value = 44369.43783564815
$ws.Cells.Item(39, 3).Value = 46070
$ws.Cells.Item(39, 7).Value = 4.8

# Row 40: A 23294-2025 -> A 63664-2023
$ws.Cells.Item(40, 1).Value = "A 63664-2023"
$ws.Cells.Item(40, 2).Value = 45275.62074074074
$ws.Cells.Item(40, 3).Value = 46070
$ws.Cells.Item(40, 7).Value = 3.2

# Row 41: A 23301-2025 -> A 57955-2024
$ws.Cells.Item(41, 1).Value = "A 57955-2024"
$ws.Cells.Item(41, 2).Value = 45631.56939814815
$ws.Cells.Item(41, 3).Value = 46070
